# qkd_simulation_inputs.xlsx update: insert a "batches" column, annotate a
# few runs with scheduling notes, and append runs 11/12 (rows 114-119).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before F ("batches"); this shifts the old
#        F:I (SKR / comment / value / Zeit) one column right to G:J and
#        fixes up existing formula references (H->I) automatically.
$ws.Columns("F").Insert()

# --- 2. Header + "batches" values for every existing data row (2-113).
$ws.Range("F1").Value2 = "batches"
$ws.Range("F2:F113").Value2 = 50

# --- 3. Scheduling / status notes added on top of existing rows.
$ws.Range("J90").Value2 = "läuft hoffentlich 6777853"
$ws.Range("L90").Value2 = "12 Uhr Donnerstag?"

$ws.Range("J102").Value2 = "läuft 6777850 "
$ws.Range("K102").Value2 = "ab 13:00"
$ws.Range("L102").Value2 = "10 Uhr Donnerstag"

# --- 4. New runs 11 and 12 (rows 114-119).
$ws.Range("A114").Value2 = 11
$ws.Range("B114").Value2 = 1
$ws.Range("C114").Value2 = -12
$ws.Range("D114").Value2 = 0.65
$ws.Range("F114").Value2 = 100
$ws.Range("I114").Value2 = 0.7
$ws.Range("E114").Formula = "=D114*I114"
$ws.Range("J114").Value2 = "kurze runden heute abend"

$ws.Range("A115").Value2 = 11
$ws.Range("B115").Value2 = 2
$ws.Range("C115").Value2 = -12
$ws.Range("D115").Value2 = 0.65
$ws.Range("F115").Value2 = 100
$ws.Range("I115").Value2 = 0.8
$ws.Range("E115").Formula = "=D115*I115"

$ws.Range("A116").Value2 = 11
$ws.Range("B116").Value2 = 3
$ws.Range("C116").Value2 = -12
$ws.Range("D116").Value2 = 0.65
$ws.Range("F116").Value2 = 100
$ws.Range("I116").Value2 = 0.9
$ws.Range("E116").Formula = "=D116*I116"

$ws.Range("A117").Value2 = 12
$ws.Range("B117").Value2 = 4
$ws.Range("C117").Value2 = -12
$ws.Range("D117").Value2 = 0.7
$ws.Range("F117").Value2 = 100
$ws.Range("I117").Value2 = 0.7
$ws.Range("E117").Formula = "=D117*I117"

$ws.Range("A118").Value2 = 12
$ws.Range("B118").Value2 = 5
$ws.Range("C118").Value2 = -12
$ws.Range("D118").Value2 = 0.7
$ws.Range("F118").Value2 = 100
$ws.Range("I118").Value2 = 0.8
$ws.Range("E118").Formula = "=D118*I118"

$ws.Range("A119").Value2 = 12
$ws.Range("B119").Value2 = 6
$ws.Range("C119").Value2 = -12
$ws.Range("D119").Value2 = 0.7
$ws.Range("F119").Value2 = 100
$ws.Range("I119").Value2 = 0.9
$ws.Range("E119").Formula = "=D119*I119"

# --- 5. Scroll / selection like the author left it (best effort: this
#        sandbox does not persist topLeftCell, only the active selection).
$ws.Activate()
$ws.Range("A106").Select()
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E116").Select()
